# Mi logt y task para semana 2 de ciclo 2
# Add two more weekly log rows (week/semana entries) to the "schedule" sheet:
#   row 3: Id=37, Semana=2, Team Leader time = 110/60 hours
#   row 4: Id=45, Semana=2, Team Leader time = 40/60 hours
# Columns A/C/E already carry per-column styles (style 1 for A/C, style 2 for E)
# so plain Range writes pick those up automatically, matching rows 1-2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 37
$ws.Range("C3").Value = 2
$ws.Range("E3").Formula = "=110/60"

# Row 4
$ws.Range("A4").Value = 45
$ws.Range("C4").Value = 2
$ws.Range("E4").Formula = "=40/60"

# Move/confirm the active selection the way the author left it (one row
# below the newly entered data).
$ws.Range("A5").Select()
